# Add 2022-Q3 data:
#  1. Insert a new worksheet named "2022-Q3" right after "总计".
#  2. Populate it with the fund holdings table for 2022-Q3.
#  3. Insert a new summary row in "总计" for 2022-Q3 (row shifts the rest down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet right after "总计" (position 2).
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Fill in the header row.
# ---------------------------------------------------------------------------
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($j = 0; $j -lt $headers.Count; $j++) {
    $q3.Cells.Item(1, $j + 2).Value = $headers[$j]
}

# ---------------------------------------------------------------------------
# Fund holdings data rows (A = index, B..H as in the source table).
# D, E, F, G are stored as text in the source workbook (e.g. "6.14"), while
# A and H are numeric. Row 16's G value is a genuine number 0 (not "0.00").
# ---------------------------------------------------------------------------
$rows = @(
    @("002666","前海开源沪港深创新成长灵活配置混合A","6.14","89.54","3.53","0.2167",10),
    @("200006","长城消费增值混合","5.46","90.90","3.60","0.1966",7),
    @("001518","万家瑞兴灵活配置混合A","2.81","83.43","3.19","0.0896",6),
    @("002667","前海开源沪港深创新成长灵活配置混合C","2.27","89.54","3.53","0.0801",10),
    @("010797","长城优选回报六个月持有期混合A","2.83","28.94","1.03","0.0291",8),
    @("005493","鑫元价值精选灵活配置混合A","0.55","76.82","3.00","0.0165",7),
    @("006241","中融医疗健康精选混合C","0.46","88.49","3.43","0.0158",4),
    @("006193","鑫元核心资产股票A","0.11","83.48","3.02","0.0033",10),
    @("010798","长城优选回报六个月持有期混合C","0.30","28.94","1.03","0.0031",8),
    @("006240","中融医疗健康精选混合A","0.05","88.49","3.43","0.0017",4),
    @("004212","中融量化智选混合A","0.03","93.01","3.05","0.0009",2),
    @("004783","中融量化智选混合C","0.01","93.01","3.05","0.0003",2),
    @("005494","鑫元价值精选灵活配置混合C","0.01","76.82","3.00","0.0003",7),
    @("006194","鑫元核心资产股票C","0.01","83.48","3.02","0.0003",10),
    @("015390","万家瑞兴灵活配置混合C","0.00","83.43","3.19",0,6)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    # A: numeric row index (0-based)
    $q3.Cells.Item($r, 1).Value = $i

    # B: fund code -- force text so leading zeros survive
    $q3.Cells.Item($r, 2).NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value = $row[0]

    # C: fund name (never numeric-looking, plain text)
    $q3.Cells.Item($r, 3).Value = $row[1]

    # D, E, F: text-formatted decimal numbers
    $q3.Cells.Item($r, 4).NumberFormat = "@"
    $q3.Cells.Item($r, 4).Value = $row[2]
    $q3.Cells.Item($r, 5).NumberFormat = "@"
    $q3.Cells.Item($r, 5).Value = $row[3]
    $q3.Cells.Item($r, 6).NumberFormat = "@"
    $q3.Cells.Item($r, 6).Value = $row[4]

    # G: text-formatted number, except row 16 which is a real 0
    if ($i -eq 14) {
        $q3.Cells.Item($r, 7).Value = $row[5]
    } else {
        $q3.Cells.Item($r, 7).NumberFormat = "@"
        $q3.Cells.Item($r, 7).Value = $row[5]
    }

    # H: numeric rank
    $q3.Cells.Item($r, 8).Value = $row[6]
}

$q3.Range("A1").Select()

# ---------------------------------------------------------------------------
# 3. Insert the 2022-Q3 summary row into "总计" (row 2, pushing others down).
# ---------------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()
$zongji.Cells.Item(2, 1).Value = 0
$zongji.Cells.Item(2, 2).Value = "2022-Q3"
$zongji.Cells.Item(2, 3).Value = 15
$zongji.Cells.Item(2, 4).Value = 0.65
